$d = $word.ActiveDocument
$p = $d.Paragraphs(1)
$r = $p.Range

# Collapse to just before the paragraph mark (end of "... para demo")
$r.End = $r.End - 1
$r.Collapse(0)

$runsToAdd = @(
    " ",
    "por que",
    " ",
    "esta",
    " muy bueno esta  en un demo de la rama ",
    "checkoup",
    " ",
    "fa",
    " "
)

foreach ($t in $runsToAdd) {
    $r.InsertAfter($t)
    $r.Collapse(0)
}

Write-Host "Full content:" $d.Content.Text
